# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (same column layout as the other
#    quarterly sheets) right before the "总计" summary sheet, and fill it
#    with the single new holding row.
# 2. Update the "总计" summary sheet with a new leading row for 2022-Q1
#    and shift the previous rows down.
#
# NOTE: worksheet object handles returned by this COM shim are resolved
# positionally. Once the sheet collection is mutated (Add/Delete/Move),
# any previously-captured Worksheet reference can silently point at the
# wrong sheet. To stay safe we always re-resolve sheets **by name** right
# before every single use via the WS helper below.

function WS([string]$name) {
    return $wb.Worksheets.Item($name)
}

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# Step 1: insert "2022-Q1" sheet right before "总计", cloning the
#         formatting of the most recent quarterly sheet ("2021-Q4").
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add((WS "总计"))
$newSheet.Name = "2022-Q1"

(WS "2021-Q4").Range("A1:H2").Copy()
(WS "2022-Q1").Range("A1:H2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------
# Scratch sheet / cell used to force numeric-looking strings (fund
# codes, percentages, …) to be stored as TEXT instead of being
# auto-coerced to numbers, mirroring how the rest of the workbook
# stores these columns. Appended at the very end so it never shifts
# the position of any sheet referenced earlier.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$scratch.Name = "__scratch__"
$scratchCell = (WS "__scratch__").Range("A1")
$scratchCell.NumberFormat = "@"

function Set-TextValue($range, [string]$text) {
    $scratchCell.Value = $text
    $scratchCell.Copy()
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---------------------------------------------------------------------
# Step 2: header row for "2022-Q1"
# ---------------------------------------------------------------------
(WS "2022-Q1").Range("B1").Value = "基金代码"
(WS "2022-Q1").Range("C1").Value = "基金名称"
(WS "2022-Q1").Range("D1").Value = "基金规模"
(WS "2022-Q1").Range("E1").Value = "股票总仓位"
(WS "2022-Q1").Range("F1").Value = "仓位占比"
(WS "2022-Q1").Range("G1").Value = "持有市值(亿元)"
(WS "2022-Q1").Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# Step 3: single holding row for "2022-Q1"
# ---------------------------------------------------------------------
(WS "2022-Q1").Range("A2").Value = 0
Set-TextValue (WS "2022-Q1").Range("B2") "512590"
(WS "2022-Q1").Range("C2").Value = "浦银安盛中证高股息精选ETF"
Set-TextValue (WS "2022-Q1").Range("D2") "0.59"
Set-TextValue (WS "2022-Q1").Range("E2") "96.43"
Set-TextValue (WS "2022-Q1").Range("F2") "2.06"
Set-TextValue (WS "2022-Q1").Range("G2") "0.0122"
(WS "2022-Q1").Range("H2").Value = 9

# ---------------------------------------------------------------------
# Step 4: update "总计" summary sheet - extend formatting to row 5 and
#         rewrite the four data rows (2022-Q1 first, followed by the
#         previously-existing quarters shifted down by one row).
# ---------------------------------------------------------------------
(WS "总计").Range("A2:D2").Copy()
(WS "总计").Range("A5:D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

(WS "总计").Range("A2").Value = 0
Set-TextValue (WS "总计").Range("B2") "2022-Q1"
(WS "总计").Range("C2").Value = 1
(WS "总计").Range("D2").Value = 0.01

(WS "总计").Range("A3").Value = 1
Set-TextValue (WS "总计").Range("B3") "2021-Q4"
(WS "总计").Range("C3").Value = 4
(WS "总计").Range("D3").Value = 0.2

(WS "总计").Range("A4").Value = 2
Set-TextValue (WS "总计").Range("B4") "2021-Q2"
(WS "总计").Range("C4").Value = 1
(WS "总计").Range("D4").Value = 0.04

(WS "总计").Range("A5").Value = 3
Set-TextValue (WS "总计").Range("B5") "2021-Q1"
(WS "总计").Range("C5").Value = 2
(WS "总计").Range("D5").Value = 0.11

# ---------------------------------------------------------------------
# cleanup
# ---------------------------------------------------------------------
(WS "__scratch__").Delete() | Out-Null
